$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-31 (columns A-H): timestamp, label, ax, ay, az, gx, gy, gz
# This represents a rolling sensor-reading window: 6 new readings were inserted at the
# top (now rows 2-7), the previously-existing 20 readings shifted down by 6 rows
# (now rows 8-27), and 4 brand-new readings were appended at the bottom (rows 28-31).
$data = @"
0,falling,-0.1190295219421371,3.933163046836853,-1.068776398897171,-0.0641408488154411,0.0114537235349416,-0.641255795955658
100,falling,1.005563497543334,4.189095079898834,-1.408742964267732,-0.0258090570569038,-0.0752891451120376,-0.038026362657547
200,falling,0.09459304809570292,4.345824718475342,-2.468320369720459,-0.2993239760398865,0.2926044464111328,0.0584903471171855
300,falling,-0.07337808609008795,4.227709949016571,-1.961481541395187,-0.08659014850854869,-0.0867428630590438,0.1705841124057769
400,falling,0.1517496109008791,4.091286063194275,-1.76620078086853,0.1209513172507286,-0.2257147133350372,-0.1018617823719978
500,falling,0.5791012048721321,4.162412345409393,-1.810639068484306,-0.070249505341053,-0.1577559560537338,0.0288633834570646
600,falling,0.4420366287231431,4.159920692443848,-1.852656066417694,0.0030543261673301,-0.0183259565383195,-0.0439822971820831
700,falling,-0.3119066953659063,3.965313136577606,-1.959281176328659,-0.024892758578062,0.007941247895359899,0.039248090237379
800,falling,-1.328240633010866,3.452776372432707,-2.262899994850159,0.0134390350431203,0.0357356183230876,0.0925460830330848
900,falling,-2.280377209186554,2.56487038731575,-2.792269647121429,-0.2785545587539673,-0.0339030213654041,0.0149661982432007
1000,falling,-1.962558269500729,2.719248294830324,-4.176476478576662,-0.2206750661134719,0.6401867866516113,0.2525927722454071
1100,falling,-0.4654391407966572,3.728965997695926,-4.883096873760225,0.9094256162643432,-0.5726861953735352,0.6157521605491638
1200,falling,-0.8071434497833259,3.663023471832275,-7.005156159400943,0.3003929853439331,-0.3223841190338135,-0.0675006061792373
1300,falling,2.102240324020393,3.889697194099427,-1.86714211106299,0.40943244099617,0.7831292152404785,-0.2806925773620605
1400,falling,0.8947885036468448,4.364453077316285,5.636359572410591,-0.4190535545349121,-0.07376197725534429,-2.578767538070679
1500,falling,2.309004843235045,3.808936268091196,8.007967770099629,2.760499954223633,-4.066987991333008,1.048855662345886
1600,falling,7.03219151496884,2.626607954502109,4.059313654899583,0.5668829679489136,-1.699732542037964,1.101084589958191
1700,falling,-0.8205822110175988,3.666950583457956,0.5117366313934346,-0.7536550164222717,-0.4547891616821289,0.0024434609804302
1800,falling,1.767575740814208,5.950162887573243,1.030251502990718,-0.0832303911447525,-0.0516181141138076,-0.2286163121461868
1900,falling,1.609463393688202,6.117009162902832,-0.02948403358459517,0.06688974797725671,0.5899431109428406,-0.3481931984424591
2000,falling,1.414643406867981,3.881679594516751,1.453149795532228,-0.7376197576522827,1.138958215713501,-0.1828014254570007
2100,falling,0.7576336860656723,3.080292761325836,0.4208878278732262,-1.0334312915802,1.02701723575592,0.1038470938801765
2200,falling,0.7153766155242933,3.394917845726014,-0.9988620877265951,-0.1331686228513717,2.704300403594971,-1.112385630607605
2300,falling,1.686916232109074,3.847963333129884,-1.479923203587532,-0.6466008424758911,-0.0087048299610614,-0.5677992701530457
2400,falling,1.892059445381164,3.739774823188782,-1.792329192161561,-1.914909839630127,-1.065807104110718,2.838385343551636
2500,falling,-0.6834564805030802,2.868067562580109,-2.777485311031344,0.1142318025231361,-0.2412917762994766,0.7263187766075134
2600,falling,-0.1287679672241191,2.86781930923462,-3.374618053436277,-0.4198171496391296,0.1876883506774902,0.732122004032135
2700,falling,0.3276743888855003,3.037657141685487,-2.661090373992916,0.0713185146450996,-0.2884811162948608,0.1009454801678657
2800,falling,0.1350092887878401,3.256921470165254,-2.639114320278168,0.0609338097274303,-0.087353728711605,0.0229074470698833
2900,falling,0.1058316230773927,3.051037281751632,-2.528261989355087,0.015118914656341,-0.0455094613134861,0.119576871395111
"@

$rows = $data -split "`n" | Where-Object { $_.Trim().Length -gt 0 }

$r = 2
foreach ($line in $rows) {
    $fields = $line -split ","
    $ws.Cells.Item($r, 1).Value = [double]$fields[0]
    $ws.Cells.Item($r, 2).Value = $fields[1]
    $ws.Cells.Item($r, 3).Value = [double]$fields[2]
    $ws.Cells.Item($r, 4).Value = [double]$fields[3]
    $ws.Cells.Item($r, 5).Value = [double]$fields[4]
    $ws.Cells.Item($r, 6).Value = [double]$fields[5]
    $ws.Cells.Item($r, 7).Value = [double]$fields[6]
    $ws.Cells.Item($r, 8).Value = [double]$fields[7]
    $r++
}
